$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - Брюн Феликс (ДЗ_1..ДЗ_4 = 5)
$ws.Range("C5:F5").Value = 5

# Row 13 - Зокирова Рохатой (ДЗ_1..ДЗ_5 = 5)
$ws.Range("C13:G13").Value = 5

# Row 16 - Кудрявцева Полина (ДЗ_1..ДЗ_5 = 5)
$ws.Range("C16:G16").Value = 5

# Row 17 - Масленникова Анастасия (ДЗ_1..ДЗ_5 = 5)
$ws.Range("C17:G17").Value = 5

# Row 26 - Теплюк Дмитрий (ДЗ_1..ДЗ_3 = 5)
$ws.Range("C26:E26").Value = 5

# Row 29 - Хабибулина Майя (ДЗ_1..ДЗ_5 = 5)
$ws.Range("C29:G29").Value = 5

# Row 31 - Юшина Полина (ДЗ_1..ДЗ_5 = 5)
$ws.Range("C31:G31").Value = 5

# Move the active selection to H5, matching the recorded sheet view state
$ws.Range("H5").Select()
